$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, pushing all existing rows (including the
# header row that was in row 1) down by one.
$ws.Rows.Item(1).Insert()

# Copy the formatting from the row that used to be the header (now row 2,
# which carries the bold/centered/bordered style) onto the new row 1.
$ws.Range("A2:L2").Copy()
$ws.Range("A1:L1").PasteSpecial(-4122)

# Populate the newly inserted row 1 with a simple numeric index sequence.
for ($col = 1; $col -le 12; $col++) {
    $ws.Cells.Item(1, $col).Value = $col - 1
}
